# Regenerate handback status report: refresh the handoff/handback timestamps
# for the second data row (4d1e0d08-ad0f-401d-a993-5f2f03cd48cf.*) across the
# Overview summary sheet and each per-locale detail sheet.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the 4d1e0d08... file.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-01 10:55:53"

# zh-cn detail sheet: "Correspond Handoff Datetime" (H3) and
# "Correspond Handback DateTime" (K3) for the 4d1e0d08... row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-01 10:55:49"
$wsZhCn.Range("K3").Value = "2016-09-01 10:56:17"

# de-de detail sheet: same two columns for the 4d1e0d08... row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-01 10:55:53"
$wsDeDe.Range("K3").Value = "2016-09-01 10:56:24"
